# Daily attendance processing - 2026-01-01 11:01:52
# Swap the order of "System" and the recorder's email address in the
# "Recorded By" column (column G) wherever it still reads
# "System, dnasr281@gmail.com", changing it to
# "dnasr281@gmail.com, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$oldText = "System, dnasr281@gmail.com"
$newText = "dnasr281@gmail.com, System"

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -eq $oldText) {
        $cell.Value2 = $newText
    }
}
